$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (H1) onto the two new
# header cells so they match the bold/centered/bordered header style,
# then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the two new columns (I0, IF), rows 2-7
$iValues = @(4, 1, 1, 1, 1, 5)
$jValues = @(8, 5, 5, 6, 7, 9)

for ($r = 2; $r -le 7; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
